$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 39 and 40 had their match data swapped (id/div/date/result-type stay put,
# but the "which match is which" columns flip between the two rows).
$row1 = 39
$row2 = 40

$cols = @("B","F","G","H","I","K","L","M","N","O","P","Q","R","S","T","U","V","W","Z","AA","AB","AC")

foreach ($col in $cols) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"

    $val1 = $ws.Range($addr1).Value2
    $val2 = $ws.Range($addr2).Value2

    $ws.Range($addr1).Value2 = $val2
    $ws.Range($addr2).Value2 = $val1
}
